$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = Get-Date -Year 2017 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B11").Value = "implémente navigation + création des différentes vues"
$ws.Range("C11").Value = 4.75

$ws.Range("C12").Select()
